# Regenerate merged AHB files
# - Rename header labels from *_old/*_new to *_FV2310/*_FV2404
# - Wrap the data range in an Excel Table (ListObject) named "Table1"
# - Freeze the header row (row 1) and set the active selection under it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header cells (row 1, columns A:K use the "old" suffix,
#    columns L:U use the "new" suffix; column K is the constant "diff").
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# 2) Turn the used range into a proper Excel Table named "Table1"
$tableRange = $ws.Range("A1:U58")
$listObject = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$listObject.Name = "Table1"

# 3) Freeze the header row and select the cell right below the freeze point
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
